$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    1  = 0.09875120888855804
    2  = -0.0059999999658515435
    3  = -0.0039999999734945391
    4  = -0.0079999999503392871
    5  = -0.0029999999777441388
    6  = -0.001999999982547962
    7  = -0.0099999999380715465
    8  = -0.0099999999384694505
    9  = -0.0019999999850552896
    10 = -0.0019999999871416207
    11 = -0.0029999999818990375
    12 = -0.0034999999798794867
    13 = 0.041668227824816384
    14 = -0.0079999999583195702
    15 = -0.0009999999970764506
    16 = -0.0019999999910371713
    17 = -0.0019999999904900534
    18 = -0.0039999999790909513
    19 = -0.003999999976889157
    20 = -0.0039999999755924165
    21 = -0.0039999999754201099
    22 = 0.0001884016048290249
    23 = -0.0049999999660261807
    24 = -0.019999999878941743
    25 = -0.019999999877229335
    26 = -0.002499999979651335
    27 = -0.0024999999795203287
    28 = -0.0019999999821340708
    29 = -0.0069999999541954239
    30 = -0.035914727056370666
    31 = -0.00287413808691106
    32 = 0.049594559904416968
    33 = -0.0039999999768625116
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
